{"js": "// Replace the date line and the 25 division-problem answers in the table\n// with their updated values, matching the target diff exactly.\nconst replacements = [\n  [\"2025-11-30 Sunday\", \"2025-12-01 Monday\"],\n  [\"66\u00f78=8, 2\", \"69\u00f74=17, 1\"],\n  [\"60\u00f79=6, 6\", \"75\u00f73=25, 0\"],\n  [\"72\u00f78=9, 0\", \"19\u00f73=6, 1\"],\n  [\"29\u00f79=3, 2\", \"46\u00f79=5, 1\"],\n  [\"74\u00f76=12, 2\", \"40\u00f78=5, 0\"],\n  [\"89\u00f78=11, 1\", \"86\u00f73=28, 2\"],\n  [\"91\u00f78=11, 3\", \"10\u00f77=1, 3\"],\n  [\"77\u00f75=15, 2\", \"30\u00f72=15, 0\"],\n  [\"53\u00f73=17, 2\", \"61\u00f74=15, 1\"],\n  [\"65\u00f73=21, 2\", \"51\u00f78=6, 3\"],\n  [\"82\u00f74=20, 2\", \"50\u00f73=16, 2\"],\n  [\"90\u00f78=11, 2\", \"61\u00f77=8, 5\"],\n  [\"59\u00f77=8, 3\", \"90\u00f73=30, 0\"],\n  [\"39\u00f73=13, 0\", \"17\u00f76=2, 5\"],\n  [\"79\u00f75=15, 4\", \"35\u00f72=17, 1\"],\n  [\"34\u00f78=4, 2\", \"23\u00f76=3, 5\"],\n  [\"73\u00f74=18, 1\", \"99\u00f74=24, 3\"],\n  [\"11\u00f73=3, 2\", \"50\u00f73=16, 2\"],\n  [\"77\u00f79=8, 5\", \"29\u00f75=5, 4\"],\n  [\"32\u00f79=3, 5\", \"51\u00f77=7, 2\"],\n  [\"80\u00f75=16, 0\", \"30\u00f77=4, 2\"],\n  [\"54\u00f76=9, 0\", \"24\u00f73=8, 0\"],\n  [\"27\u00f73=9, 0\", \"78\u00f77=11, 1\"],\n  [\"28\u00f74=7, 0\", \"79\u00f76=13, 1\"],\n  [\"83\u00f79=9, 2\", \"41\u00f73=13, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 division-problem answers to the\n# new values, matching the target diff exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2025-11-30 Sunday\"; Replace = \"2025-12-01 Monday\" },\n    @{ Find = \"66\u00f78=8, 2\";  Replace = \"69\u00f74=17, 1\" },\n    @{ Find = \"60\u00f79=6, 6\";  Replace = \"75\u00f73=25, 0\" },\n    @{ Find = \"72\u00f78=9, 0\";  Replace = \"19\u00f73=6, 1\" },\n    @{ Find = \"29\u00f79=3, 2\";  Replace = \"46\u00f79=5, 1\" },\n    @{ Find = \"74\u00f76=12, 2\"; Replace = \"40\u00f78=5, 0\" },\n    @{ Find = \"89\u00f78=11, 1\"; Replace = \"86\u00f73=28, 2\" },\n    @{ Find = \"91\u00f78=11, 3\"; Replace = \"10\u00f77=1, 3\" },\n    @{ Find = \"77\u00f75=15, 2\"; Replace = \"30\u00f72=15, 0\" },\n    @{ Find = \"53\u00f73=17, 2\"; Replace = \"61\u00f74=15, 1\" },\n    @{ Find = \"65\u00f73=21, 2\"; Replace = \"51\u00f78=6, 3\" },\n    @{ Find = \"82\u00f74=20, 2\"; Replace = \"50\u00f73=16, 2\" },\n    @{ Find = \"90\u00f78=11, 2\"; Replace = \"61\u00f77=8, 5\" },\n    @{ Find = \"59\u00f77=8, 3\";  Replace = \"90\u00f73=30, 0\" },\n    @{ Find = \"39\u00f73=13, 0\"; Replace = \"17\u00f76=2, 5\" },\n    @{ Find = \"79\u00f75=15, 4\"; Replace = \"35\u00f72=17, 1\" },\n    @{ Find = \"34\u00f78=4, 2\";  Replace = \"23\u00f76=3, 5\" },\n    @{ Find = \"73\u00f74=18, 1\"; Replace = \"99\u00f74=24, 3\" },\n    @{ Find = \"11\u00f73=3, 2\";  Replace = \"50\u00f73=16, 2\" },\n    @{ Find = \"77\u00f79=8, 5\";  Replace = \"29\u00f75=5, 4\" },\n    @{ Find = \"32\u00f79=3, 5\";  Replace = \"51\u00f77=7, 2\" },\n    @{ Find = \"80\u00f75=16, 0\"; Replace = \"30\u00f77=4, 2\" },\n    @{ Find = \"54\u00f76=9, 0\";  Replace = \"24\u00f73=8, 0\" },\n    @{ Find = \"27\u00f73=9, 0\";  Replace = \"78\u00f77=11, 1\" },\n    @{ Find = \"28\u00f74=7, 0\";  Replace = \"79\u00f76=13, 1\" },\n    @{ Find = \"83\u00f79=9, 2\";  Replace = \"41\u00f73=13, 2\" }\n)\n\nforeach ($item in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    [void]$rng.Find.Execute($item.Find, $false, $false, $false, $false, $false, $true, 1, $false, $item.Replace, 2)\n}\n"}
